$d = $word.ActiveDocument

# Move to the very end of the document body (collapsed range at the end)
$endRange = $d.Content
$endRange.Collapse(0)  # wdCollapseEnd

for ($i = 0; $i -lt 7; $i++) {
    $p = $endRange.Paragraphs.Add($endRange)
    $p.Range.Font.Size = 22
}
